# Daily attendance processing - 2025-10-10 14:46:29
# Re-orders the "Recorded By" (column G) contributor lists for a set of rows
# on the "Session Analysis Results" sheet, without altering any other data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$updates = @{
    'G2'   = 'system, backup@backdoor.com, System'
    'G3'   = 'dnasr281@gmail.com, System'
    'G6'   = 'dnasr281@gmail.com, System'
    'G7'   = 'admin@admin.com, System'
    'G10'  = 'dnasr281@gmail.com, System'
    'G11'  = 'dnasr281@gmail.com, System'
    'G12'  = 'dnasr281@gmail.com, System'
    'G13'  = 'dnasr281@gmail.com, System'
    'G14'  = 'dnasr281@gmail.com, System'
    'G15'  = 'dnasr281@gmail.com, System'
    'G29'  = 'system, backup@backdoor.com, System'
    'G30'  = 'dnasr281@gmail.com, System'
    'G33'  = 'dnasr281@gmail.com, System'
    'G34'  = 'admin@admin.com, System'
    'G37'  = 'dnasr281@gmail.com, System'
    'G38'  = 'dnasr281@gmail.com, System'
    'G39'  = 'dnasr281@gmail.com, System'
    'G40'  = 'dnasr281@gmail.com, System'
    'G41'  = 'dnasr281@gmail.com, System'
    'G42'  = 'dnasr281@gmail.com, System'
    'G56'  = 'system, backup@backdoor.com, System'
    'G57'  = 'dnasr281@gmail.com, System'
    'G60'  = 'dnasr281@gmail.com, System'
    'G61'  = 'admin@admin.com, System'
    'G64'  = 'dnasr281@gmail.com, System'
    'G65'  = 'dnasr281@gmail.com, System'
    'G66'  = 'dnasr281@gmail.com, System'
    'G67'  = 'dnasr281@gmail.com, System'
    'G68'  = 'dnasr281@gmail.com, System'
    'G69'  = 'dnasr281@gmail.com, System'
    'G86'  = 'dnasr281@gmail.com, System'
    'G87'  = 'dnasr281@gmail.com, System'
    'G88'  = 'dnasr281@gmail.com, System'
    'G89'  = 'dnasr281@gmail.com, System'
    'G90'  = 'dnasr281@gmail.com, admin@admin.com'
    'G93'  = 'dnasr281@gmail.com, System'
    'G95'  = 'dnasr281@gmail.com, System'
    'G112' = 'dnasr281@gmail.com, System'
    'G113' = 'dnasr281@gmail.com, System'
    'G114' = 'dnasr281@gmail.com, System'
    'G115' = 'dnasr281@gmail.com, System'
    'G116' = 'dnasr281@gmail.com, admin@admin.com'
    'G119' = 'dnasr281@gmail.com, System'
    'G121' = 'dnasr281@gmail.com, System'
    'G138' = 'dnasr281@gmail.com, System'
    'G139' = 'dnasr281@gmail.com, System'
    'G140' = 'dnasr281@gmail.com, System'
    'G141' = 'dnasr281@gmail.com, System'
    'G142' = 'dnasr281@gmail.com, admin@admin.com'
    'G145' = 'dnasr281@gmail.com, System'
    'G147' = 'dnasr281@gmail.com, System'
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
